# Update countries & provincias Spain
# - Swap shared-string order of "Serbia" / "Filipinas" (country names in rows 43/44
#   are swapped so that row 43 now shows Serbia and row 44 shows Filipinas)
# - Update the "Datos actualizados ..." timestamp string
# - Refresh several countries' statistics (columns B..H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 15:33"

# Swap country names for rows 43 and 44 (Filipinas <-> Serbia)
$ws.Range("A43").Value = "Serbia"
$ws.Range("A44").Value = "Filipinas"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1189845
$ws.Range("C4").Value = 1723
$ws.Range("E4").Value = 942541
$ws.Range("G4").Value = 35
$ws.Range("H4").Value = 68633

# Espana (row 5)
$ws.Range("B5").Value = 248301
$ws.Range("C5").Value = 1179
$ws.Range("D5").Value = 151633
$ws.Range("E5").Value = 71240
$ws.Range("F5").Value = 2254
$ws.Range("G5").Value = 164
$ws.Range("H5").Value = 25428

# Row 43 (now Serbia) - updated stats
$ws.Range("B43").Value = 9557
$ws.Range("C43").Value = 93
$ws.Range("D43").Value = 1574
$ws.Range("E43").Value = 7786
$ws.Range("F43").Value = 54
$ws.Range("G43").Value = 4
$ws.Range("H43").Value = 197

# Row 44 (now Filipinas) - updated stats
$ws.Range("B44").Value = 9485
$ws.Range("C44").Value = 262
$ws.Range("D44").Value = 1315
$ws.Range("E44").Value = 7547
$ws.Range("F44").Value = 31
$ws.Range("G44").Value = 16
$ws.Range("H44").Value = 623

# Noruega (row 46)
$ws.Range("B46").Value = 7884
$ws.Range("C46").Value = 37
$ws.Range("E46").Value = 7638
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 214

# Sri Lanka (row 103)
$ws.Range("B103").Value = 719
$ws.Range("C103").Value = 1
$ws.Range("E103").Value = 525

# Mauricio (row 127)
$ws.Range("D127").Value = 316
$ws.Range("E127").Value = 6

# Uganda (row 158)
$ws.Range("D158").Value = 55
$ws.Range("E158").Value = 34
